# feat: add 2022-Q3 data
#
# 1) Duplicate the existing "2021-Q1" sheet (keeps its header/row styling
#    intact) into a brand-new "2022-Q3" sheet positioned right after "总计".
# 2) Overwrite the duplicate's fund-metrics with the new quarter's numbers.
# 3) Insert a new summary row on "总计" for 2022-Q3, pushing the older
#    quarters down by one row.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)          # 总计
$q1sheet = $wb.Worksheets.Item("2021-Q1")  # template for formatting/layout

# --- 1) new "2022-Q3" sheet, cloned from "2021-Q1" so it inherits the same
#        header/border/alignment styling, then dropped right after it (i.e.
#        right after 总计, before the original 2021-Q1 sheet). ---
$q1sheet.Copy($q1sheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# --- 2) update the new sheet's figures for 2022-Q3 ---
$newSheet.Range("D1").Value = "基金规模"

$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.12"
$newSheet.Range("E2").Value = "90.06"
$newSheet.Range("F2").Value = "2.70"
$newSheet.Range("G2").Value = "0.0302"
$newSheet.Range("H2").Value = 8

# --- 3) insert the matching row into 总计, shifting 2021-Q1/2020-Q4 down ---
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.03

# the index column is a plain 0-based counter, not a formula, so renumber
# the rows that got pushed down by the insert above
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
